# Auto-optimize exam scheduling: dynamically adjusts exams per slot (1-4)
# to guarantee all courses are scheduled within date range.
#
# Applies the re-shuffled course placements to the Section_A / Section_B
# timetables, and updates the CS151 elective lecture slot info on the
# Elective_Coordination sheet accordingly.

$wb = $excel.ActiveWorkbook

# ---- Section_A ----
$wsA = $wb.Worksheets.Item("Section_A")

# Row 2 (09:00-10:30)
$wsA.Range("B2").Value = "DS161"
$wsA.Range("C2").Value = "Free"
$wsA.Range("D2").Value = "MA161"
$wsA.Range("E2").Value = "DS161"
$wsA.Range("F2").Value = "MA162"

# Row 3 (10:30-12:00)
$wsA.Range("B3").Value = "MA161"
$wsA.Range("C3").Value = "CS161"
$wsA.Range("D3").Value = "EC161"
$wsA.Range("E3").Value = "Free"
$wsA.Range("F3").Value = "CS151 (Elective)"

# Row 5 (13:00-14:30)
$wsA.Range("B5").Value = "HS161"
$wsA.Range("C5").Value = "EC161"
$wsA.Range("D5").Value = "MA162"
$wsA.Range("E5").Value = "Free"
$wsA.Range("F5").Value = "CS161"

# Row 7 (15:30-17:00)
$wsA.Range("B7").Value = "EC161"
$wsA.Range("C7").Value = "CS151 (Elective)"
$wsA.Range("D7").Value = "CS161"
$wsA.Range("E7").Value = "HS161"
$wsA.Range("F7").Value = "HS161"

# ---- Section_B ----
$wsB = $wb.Worksheets.Item("Section_B")

# Row 2 (09:00-10:30)
$wsB.Range("B2").Value = "EC161"
$wsB.Range("C2").Value = "EC161"
$wsB.Range("D2").Value = "MA161"
$wsB.Range("E2").Value = "HS161"
$wsB.Range("F2").Value = "CS161"

# Row 3 (10:30-12:00)
$wsB.Range("B3").Value = "Free"
$wsB.Range("C3").Value = "Free"
$wsB.Range("D3").Value = "CS161"
$wsB.Range("E3").Value = "MA161"
$wsB.Range("F3").Value = "CS151 (Elective)"

# Row 5 (13:00-14:30)
$wsB.Range("B5").Value = "Free"
$wsB.Range("C5").Value = "MA162"
$wsB.Range("D5").Value = "HS161"
$wsB.Range("E5").Value = "DS161"
$wsB.Range("F5").Value = "HS161"

# Row 7 (15:30-17:00)
$wsB.Range("B7").Value = "DS161"
$wsB.Range("C7").Value = "CS151 (Elective)"
$wsB.Range("D7").Value = "EC161"
$wsB.Range("E7").Value = "CS161"
$wsB.Range("F7").Value = "MA162"

# ---- Elective_Coordination ----
$wsE = $wb.Worksheets.Item("Elective_Coordination")

# CS151 Lecture 1 moved from Mon 09:00-10:30 to Tue 15:30-17:00
$wsE.Range("C10").Value = "Tue"
$wsE.Range("D10").Value = "15:30-17:00"

# CS151 Lecture 2 moved from Thu 15:30-17:00 to Fri 10:30-12:00
$wsE.Range("C11").Value = "Fri"
$wsE.Range("D11").Value = "10:30-12:00"
